$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Change Log")

# Row 12: set the date, update the Changes/Notes text, and grow the row to fit
# the new (longer) Notes text.
$ws.Range("A12").Value = 45941

$ws.Range("D12").Value = "Notes`n- Likely could've had ChatGPT do this as you did have it get the instruction encodings in the 1st place. But it failed when you attempted to ask it to verify them, so I verified them myself.`nBugs`n- "

$ws.Range("B12").Value = "Changes`n- MODIFIED: MyMIF.mif`n- ADDED: Verified encodings and added verification comments                                                                                                                                                     "

$ws.Rows.Item(12).RowHeight = 100.8

# Move the selection to where the author left off editing.
$ws.Range("B13").Select() | Out-Null
